$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "10+1=11"
$t.Cell(1, 2).Range.Text = "22+29=51"
$t.Cell(1, 3).Range.Text = "69-38=31"
$t.Cell(1, 4).Range.Text = "4+78=82"
$t.Cell(1, 5).Range.Text = "96-84=12"
$t.Cell(2, 1).Range.Text = "54-32=22"
$t.Cell(2, 2).Range.Text = "50-19=31"
$t.Cell(2, 3).Range.Text = "87-74=13"
$t.Cell(2, 4).Range.Text = "72+10=82"
$t.Cell(2, 5).Range.Text = "37-17=20"
$t.Cell(3, 1).Range.Text = "65+24=89"
$t.Cell(3, 2).Range.Text = "67-41=26"
$t.Cell(3, 3).Range.Text = "35-7=28"
$t.Cell(3, 4).Range.Text = "98-6=92"
$t.Cell(3, 5).Range.Text = "45-6=39"
$t.Cell(4, 1).Range.Text = "97-73=24"
$t.Cell(4, 2).Range.Text = "33+62=95"
$t.Cell(4, 3).Range.Text = "48-25=23"
$t.Cell(4, 4).Range.Text = "43-19=24"
$t.Cell(4, 5).Range.Text = "19+78=97"
$t.Cell(5, 1).Range.Text = "41-7=34"
$t.Cell(5, 2).Range.Text = "93-16=77"
$t.Cell(5, 3).Range.Text = "55+19=74"
$t.Cell(5, 4).Range.Text = "28+45=73"
$t.Cell(5, 5).Range.Text = "23+44=67"
$t.Cell(6, 1).Range.Text = "86-58=28"
$t.Cell(6, 2).Range.Text = "44-21=23"
$t.Cell(6, 3).Range.Text = "42-36=6"
$t.Cell(6, 4).Range.Text = "10-0=10"
$t.Cell(6, 5).Range.Text = "39+15=54"
$t.Cell(7, 1).Range.Text = "95-71=24"
$t.Cell(7, 2).Range.Text = "87-66=21"
$t.Cell(7, 3).Range.Text = "53-30=23"
$t.Cell(7, 4).Range.Text = "35-25=10"
$t.Cell(7, 5).Range.Text = "57+33=90"
$t.Cell(8, 1).Range.Text = "67-10=57"
$t.Cell(8, 2).Range.Text = "44+1=45"
$t.Cell(8, 3).Range.Text = "27+41=68"
$t.Cell(8, 4).Range.Text = "89-78=11"
$t.Cell(8, 5).Range.Text = "21+67=88"
$t.Cell(9, 1).Range.Text = "23+19=42"
$t.Cell(9, 2).Range.Text = "75-26=49"
$t.Cell(9, 3).Range.Text = "10+13=23"
$t.Cell(9, 4).Range.Text = "15+17=32"
$t.Cell(9, 5).Range.Text = "34-29=5"
$t.Cell(10, 1).Range.Text = "78-66=12"
$t.Cell(10, 2).Range.Text = "36+15=51"
$t.Cell(10, 3).Range.Text = "12+27=39"
$t.Cell(10, 4).Range.Text = "76-9=67"
$t.Cell(10, 5).Range.Text = "37-4=33"
$t.Cell(11, 1).Range.Text = "44-30=14"
$t.Cell(11, 2).Range.Text = "50+32=82"
$t.Cell(11, 3).Range.Text = "83-60=23"
$t.Cell(11, 4).Range.Text = "87-63=24"
$t.Cell(11, 5).Range.Text = "39-18=21"
$t.Cell(12, 1).Range.Text = "52-12=40"
$t.Cell(12, 2).Range.Text = "81-77=4"
$t.Cell(12, 3).Range.Text = "23-15=8"
$t.Cell(12, 4).Range.Text = "67-35=32"
$t.Cell(12, 5).Range.Text = "38-19=19"
$t.Cell(13, 1).Range.Text = "70-55=15"
$t.Cell(13, 2).Range.Text = "37-32=5"
$t.Cell(13, 3).Range.Text = "29+58=87"
$t.Cell(13, 4).Range.Text = "94-24=70"
$t.Cell(13, 5).Range.Text = "65+16=81"
$t.Cell(14, 1).Range.Text = "85-16=69"
$t.Cell(14, 2).Range.Text = "64+26=90"
$t.Cell(14, 3).Range.Text = "75-38=37"
$t.Cell(14, 4).Range.Text = "37+49=86"
$t.Cell(14, 5).Range.Text = "75+15=90"
$t.Cell(15, 1).Range.Text = "53-32=21"
$t.Cell(15, 2).Range.Text = "15+76=91"
$t.Cell(15, 3).Range.Text = "47+15=62"
$t.Cell(15, 4).Range.Text = "36-25=11"
$t.Cell(15, 5).Range.Text = "82-35=47"
$t.Cell(16, 1).Range.Text = "37+38=75"
$t.Cell(16, 2).Range.Text = "28+45=73"
$t.Cell(16, 3).Range.Text = "10+23=33"
$t.Cell(16, 4).Range.Text = "73-6=67"
$t.Cell(16, 5).Range.Text = "97-23=74"
$t.Cell(17, 1).Range.Text = "91-4=87"
$t.Cell(17, 2).Range.Text = "64-4=60"
$t.Cell(17, 3).Range.Text = "58-3=55"
$t.Cell(17, 4).Range.Text = "88-23=65"
$t.Cell(17, 5).Range.Text = "24+35=59"
$t.Cell(18, 1).Range.Text = "95-44=51"
$t.Cell(18, 2).Range.Text = "47+51=98"
$t.Cell(18, 3).Range.Text = "46-38=8"
$t.Cell(18, 4).Range.Text = "36+55=91"
$t.Cell(18, 5).Range.Text = "62-11=51"
$t.Cell(19, 1).Range.Text = "4+13=17"
$t.Cell(19, 2).Range.Text = "10+35=45"
$t.Cell(19, 3).Range.Text = "33+8=41"
$t.Cell(19, 4).Range.Text = "87-12=75"
$t.Cell(19, 5).Range.Text = "34-22=12"
$t.Cell(20, 1).Range.Text = "97-71=26"
$t.Cell(20, 2).Range.Text = "58+16=74"
$t.Cell(20, 3).Range.Text = "93+0=93"
$t.Cell(20, 4).Range.Text = "94-73=21"
$t.Cell(20, 5).Range.Text = "39-25=14"
